{"js": "// Replace the date line and the 25 \"three-digit \u00f7 one-digit\" table answers\n// with their new values, as a set of unique find-and-replace operations.\nconst replacements = [\n  [\"2024-04-17 Wednesday\", \"2024-04-18 Thursday\"],\n  [\"261\u00f75=52, 1\", \"804\u00f78=100, 4\"],\n  [\"225\u00f73=75, 0\", \"684\u00f76=114, 0\"],\n  [\"110\u00f76=18, 2\", \"655\u00f73=218, 1\"],\n  [\"186\u00f76=31, 0\", \"783\u00f73=261, 0\"],\n  [\"902\u00f76=150, 2\", \"295\u00f78=36, 7\"],\n  [\"299\u00f73=99, 2\", \"433\u00f78=54, 1\"],\n  [\"788\u00f76=131, 2\", \"629\u00f79=69, 8\"],\n  [\"895\u00f72=447, 1\", \"194\u00f79=21, 5\"],\n  [\"989\u00f75=197, 4\", \"251\u00f77=35, 6\"],\n  [\"933\u00f78=116, 5\", \"705\u00f77=100, 5\"],\n  [\"586\u00f79=65, 1\", \"453\u00f76=75, 3\"],\n  [\"424\u00f72=212, 0\", \"971\u00f73=323, 2\"],\n  [\"850\u00f76=141, 4\", \"114\u00f74=28, 2\"],\n  [\"917\u00f78=114, 5\", \"373\u00f74=93, 1\"],\n  [\"490\u00f73=163, 1\", \"839\u00f76=139, 5\"],\n  [\"475\u00f77=67, 6\", \"468\u00f74=117, 0\"],\n  [\"642\u00f79=71, 3\", \"491\u00f72=245, 1\"],\n  [\"787\u00f76=131, 1\", \"249\u00f73=83, 0\"],\n  [\"427\u00f76=71, 1\", \"312\u00f78=39, 0\"],\n  [\"260\u00f77=37, 1\", \"316\u00f72=158, 0\"],\n  [\"643\u00f78=80, 3\", \"239\u00f74=59, 3\"],\n  [\"491\u00f75=98, 1\", \"420\u00f74=105, 0\"],\n  [\"831\u00f74=207, 3\", \"707\u00f73=235, 2\"],\n  [\"722\u00f78=90, 2\", \"641\u00f75=128, 1\"],\n  [\"878\u00f79=97, 5\", \"914\u00f76=152, 2\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  for (const [oldText, newText] of replacements) {\n    const results = paragraph.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n    if (results.items.length > 0) {\n      for (const range of results.items) {\n        range.insertText(newText, Word.InsertLocation.replace);\n      }\n      await context.sync();\n    }\n  }\n}\n", "ps1": "# Replace the date line and the 25 \"three-digit \u00f7 one-digit\" table answers\n# with their new values, using Find/Replace on the whole document content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-17 Wednesday\", \"2024-04-18 Thursday\"),\n    @(\"261\u00f75=52, 1\", \"804\u00f78=100, 4\"),\n    @(\"225\u00f73=75, 0\", \"684\u00f76=114, 0\"),\n    @(\"110\u00f76=18, 2\", \"655\u00f73=218, 1\"),\n    @(\"186\u00f76=31, 0\", \"783\u00f73=261, 0\"),\n    @(\"902\u00f76=150, 2\", \"295\u00f78=36, 7\"),\n    @(\"299\u00f73=99, 2\", \"433\u00f78=54, 1\"),\n    @(\"788\u00f76=131, 2\", \"629\u00f79=69, 8\"),\n    @(\"895\u00f72=447, 1\", \"194\u00f79=21, 5\"),\n    @(\"989\u00f75=197, 4\", \"251\u00f77=35, 6\"),\n    @(\"933\u00f78=116, 5\", \"705\u00f77=100, 5\"),\n    @(\"586\u00f79=65, 1\", \"453\u00f76=75, 3\"),\n    @(\"424\u00f72=212, 0\", \"971\u00f73=323, 2\"),\n    @(\"850\u00f76=141, 4\", \"114\u00f74=28, 2\"),\n    @(\"917\u00f78=114, 5\", \"373\u00f74=93, 1\"),\n    @(\"490\u00f73=163, 1\", \"839\u00f76=139, 5\"),\n    @(\"475\u00f77=67, 6\", \"468\u00f74=117, 0\"),\n    @(\"642\u00f79=71, 3\", \"491\u00f72=245, 1\"),\n    @(\"787\u00f76=131, 1\", \"249\u00f73=83, 0\"),\n    @(\"427\u00f76=71, 1\", \"312\u00f78=39, 0\"),\n    @(\"260\u00f77=37, 1\", \"316\u00f72=158, 0\"),\n    @(\"643\u00f78=80, 3\", \"239\u00f74=59, 3\"),\n    @(\"491\u00f75=98, 1\", \"420\u00f74=105, 0\"),\n    @(\"831\u00f74=207, 3\", \"707\u00f73=235, 2\"),\n    @(\"722\u00f78=90, 2\", \"641\u00f75=128, 1\"),\n    @(\"878\u00f79=97, 5\", \"914\u00f76=152, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
